$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.609586333333333
$ws.Range("H2").Value = 4.828759
$ws.Range("I2").Value = 0.05107819292772156
$ws.Range("J2").Value = 0.05107819292772156
$ws.Range("M2").Value = 14.11187666666667
$ws.Range("N2").Value = 42.33562999999999
$ws.Range("O2").Value = 0.08862966207485527
$ws.Range("P2").Value = 0.08862966207485526
$ws.Range("Q2").Value = 22.71428382035222
$ws.Range("R2").Value = 204.42855438317
$ws.Range("S2").Value = 0.004527042978578225
$ws.Range("T2").Value = 0.004527042978578224

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.609586333333333
$ws.Range("H3").Value = 4.828759
$ws.Range("I3").Value = 0.05107819292772156
$ws.Range("J3").Value = 0.05107819292772156
$ws.Range("O3").Value = 0.7176943460983047
$ws.Range("P3").Value = 0.7176943460983046
$ws.Range("Q3").Value = 183.9329259742708
$ws.Range("R3").Value = 1655.396333768437
$ws.Range("S3").Value = 0.03665853027314418
$ws.Range("T3").Value = 0.03665853027314418

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.609586333333333
$ws.Range("H4").Value = 4.828759
$ws.Range("I4").Value = 0.05107819292772156
$ws.Range("J4").Value = 0.05107819292772156
$ws.Range("O4").Value = 0.1936759918268401
$ws.Range("P4").Value = 0.1936759918268401
$ws.Range("Q4").Value = 49.63588198979656
$ws.Range("R4").Value = 446.722937908169
$ws.Range("S4").Value = 0.009892619675999165
$ws.Range("T4").Value = 0.009892619675999163

# Row 5
$ws.Range("I5").Value = 0.5992082897496871
$ws.Range("J5").Value = 0.5992082897496871
$ws.Range("M5").Value = 14.11187666666667
$ws.Range("N5").Value = 42.33562999999999
$ws.Range("O5").Value = 0.08862966207485527
$ws.Range("P5").Value = 0.08862966207485526
$ws.Range("Q5").Value = 266.4657142460378
$ws.Range("R5").Value = 2398.19142821434
$ws.Range("S5").Value = 0.05310762823296673
$ws.Range("T5").Value = 0.05310762823296673

# Row 6
$ws.Range("I6").Value = 0.5992082897496871
$ws.Range("J6").Value = 0.5992082897496871
$ws.Range("O6").Value = 0.7176943460983047
$ws.Range("P6").Value = 0.7176943460983046
$ws.Range("S6").Value = 0.4300484016885852
$ws.Range("T6").Value = 0.4300484016885851

# Row 7
$ws.Range("I7").Value = 0.5992082897496871
$ws.Range("J7").Value = 0.5992082897496871
$ws.Range("O7").Value = 0.1936759918268401
$ws.Range("P7").Value = 0.1936759918268401
$ws.Range("S7").Value = 0.1160522598281353
$ws.Range("T7").Value = 0.1160522598281352

# Row 8
$ws.Range("H8").Value = 33.060729
$ws.Range("I8").Value = 0.3497135173225914
$ws.Range("J8").Value = 0.3497135173225914
$ws.Range("M8").Value = 14.11187666666667
$ws.Range("N8").Value = 42.33562999999999
$ws.Range("O8").Value = 0.08862966207485527
$ws.Range("P8").Value = 0.08862966207485526
$ws.Range("Q8").Value = 155.5163100526967
$ws.Range("R8").Value = 1399.64679047427
$ws.Range("S8").Value = 0.03099499086331032
$ws.Range("T8").Value = 0.03099499086331032

# Row 9
$ws.Range("H9").Value = 33.060729
$ws.Range("I9").Value = 0.3497135173225914
$ws.Range("J9").Value = 0.3497135173225914
$ws.Range("O9").Value = 0.7176943460983047
$ws.Range("P9").Value = 0.7176943460983046
$ws.Range("Q9").Value = 1259.320794393017
$ws.Range("S9").Value = 0.2509874141365754
$ws.Range("T9").Value = 0.2509874141365754

# Row 10
$ws.Range("H10").Value = 33.060729
$ws.Range("I10").Value = 0.3497135173225914
$ws.Range("J10").Value = 0.3497135173225914
$ws.Range("O10").Value = 0.1936759918268401
$ws.Range("P10").Value = 0.1936759918268401
$ws.Range("S10").Value = 0.06773111232270573
$ws.Range("T10").Value = 0.06773111232270573
